# BlackJack V.1.0 official submission: grow the GUI mockup sheet with a new
# row 12 placeholder (a single-space text field) and give row 3 some breathing
# room (custom height), then leave the selection parked on E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 grows a custom height of 42.75pt.
$ws.Rows.Item(3).RowHeight = 42.75

# New row 12: B12 holds a single space " " (a placeholder text field value).
$ws.Cells.Item(12, 2).Value = " "

# Move/save the active selection to E4.
$ws.Range("E4").Select()
